# Scheduled runner update: refresh market-price-derived Leve profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) across the
# per-job sheets, per the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2152.8333
$ws.Range("J28").Value = 3478
$ws.Range("L28").Value = 3478
$ws.Range("N28").Value = -4448

$ws.Range("H33").Value = 110.64286
$ws.Range("I33").Value = 94.8
$ws.Range("K33").Value = 94.8
$ws.Range("M33").Value = 134.2

$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -3350

$ws.Range("H112").Value = 6999.375
$ws.Range("J112").Value = 7842.143
$ws.Range("L112").Value = 23526.429
$ws.Range("N112").Value = -25742.429

$ws.Range("H129").Value = 833.9123
$ws.Range("I129").Value = 596.5714
$ws.Range("J129").Value = 867.14
$ws.Range("K129").Value = 1789.7142
$ws.Range("L129").Value = 2601.42
$ws.Range("M129").Value = 3210.2858
$ws.Range("N129").Value = -12601.42

$ws.Range("H132").Value = 1068.6666
$ws.Range("J132").Value = 1231.25
$ws.Range("L132").Value = 3693.75
$ws.Range("N132").Value = -8753.75

$ws.Range("H135").Value = 503.6875
$ws.Range("J135").Value = 490
$ws.Range("L135").Value = 4410
$ws.Range("N135").Value = -9480

$ws.Range("H137").Value = 1895.4736
$ws.Range("I137").Value = 1605.5454
$ws.Range("K137").Value = 4816.6362
$ws.Range("M137").Value = -2266.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3818.2856
$ws.Range("I61").Value = 2812.4
$ws.Range("J61").Value = 4978.923
$ws.Range("K61").Value = 2812.4
$ws.Range("L61").Value = 4978.923
$ws.Range("M61").Value = -2600.4
$ws.Range("N61").Value = -5402.923

$ws.Range("H102").Value = 1465.7778
$ws.Range("I102").Value = 1211.5
$ws.Range("K102").Value = 1211.5
$ws.Range("M102").Value = 410.5

$ws.Range("H110").Value = 1466.4231
$ws.Range("I110").Value = 1124.56
$ws.Range("K110").Value = 1124.56
$ws.Range("M110").Value = 920.4400000000001

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null

$ws.Range("H132").Value = 1692.7
$ws.Range("I132").Value = 1149.3334
$ws.Range("K132").Value = 3448.0002
$ws.Range("M132").Value = -918.0001999999999

$ws.Range("H136").Value = 3818.2856
$ws.Range("I136").Value = 2812.4
$ws.Range("J136").Value = 4978.923
$ws.Range("K136").Value = 8437.200000000001
$ws.Range("L136").Value = 14936.769
$ws.Range("M136").Value = -5887.200000000001
$ws.Range("N136").Value = -20036.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1103.1666
$ws.Range("I107").Value = 1103.1666
$ws.Range("K107").Value = 1103.1666
$ws.Range("M107").Value = 816.8334

$ws.Range("H134").Value = 3932.2666
$ws.Range("I134").Value = 3932.2666
$ws.Range("K134").Value = 11796.7998
$ws.Range("M134").Value = -9261.799800000001

$ws.Range("H135").Value = 52514.75
$ws.Range("J135").Value = 52514.75
$ws.Range("L135").Value = 52514.75
$ws.Range("N135").Value = -62654.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 360
$ws.Range("I7").Value = 166.66667
$ws.Range("J7").Value = 650
$ws.Range("K7").Value = 166.66667
$ws.Range("L7").Value = 650
$ws.Range("M7").Value = -53.66667000000001
$ws.Range("N7").Value = -876

$ws.Range("H16").Value = 1449.8334
$ws.Range("I16").Value = 1066.6666
$ws.Range("J16").Value = 1833
$ws.Range("K16").Value = 1066.6666
$ws.Range("L16").Value = 1833
$ws.Range("M16").Value = -779.6666
$ws.Range("N16").Value = -2407

$ws.Range("H31").Value = 2141.4243
$ws.Range("I31").Value = 1730.3529
$ws.Range("J31").Value = 2578.1875
$ws.Range("K31").Value = 1730.3529
$ws.Range("L31").Value = 2578.1875
$ws.Range("M31").Value = -1435.3529
$ws.Range("N31").Value = -3168.1875

$ws.Range("H34").Value = 2141.4243
$ws.Range("I34").Value = 1730.3529
$ws.Range("J34").Value = 2578.1875
$ws.Range("K34").Value = 1730.3529
$ws.Range("L34").Value = 2578.1875
$ws.Range("M34").Value = -1528.3529
$ws.Range("N34").Value = -2982.1875

$ws.Range("H53").Value = 68000
$ws.Range("J53").Value = 68000
$ws.Range("L53").Value = 68000
$ws.Range("N53").Value = -69214

$ws.Range("H58").Value = 1978129.4
$ws.Range("I58").Value = 2899797.2
$ws.Range("J58").Value = 3126.7144
$ws.Range("K58").Value = 2899797.2
$ws.Range("L58").Value = 3126.7144
$ws.Range("M58").Value = -2899594.2
$ws.Range("N58").Value = -3532.7144

$ws.Range("H105").Value = 2666.6667
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws.Range("H113").Value = 1449.8334
$ws.Range("I113").Value = 1066.6666
$ws.Range("J113").Value = 1833
$ws.Range("K113").Value = 1066.6666
$ws.Range("L113").Value = 1833
$ws.Range("M113").Value = 1103.3334
$ws.Range("N113").Value = -6173

$ws.Range("H122").Value = 9507
$ws.Range("I122").Value = 9000
$ws.Range("K122").Value = 27000
$ws.Range("M122").Value = -24550

$ws.Range("H136").Value = 1978129.4
$ws.Range("I136").Value = 2899797.2
$ws.Range("J136").Value = 3126.7144
$ws.Range("K136").Value = 8699391.600000001
$ws.Range("L136").Value = 9380.143199999999
$ws.Range("M136").Value = -8696841.600000001
$ws.Range("N136").Value = -14480.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2831.5
$ws.Range("I64").Value = 1744.5
$ws.Range("J64").Value = 3375
$ws.Range("K64").Value = 5233.5
$ws.Range("L64").Value = 10125
$ws.Range("M64").Value = -4963.5
$ws.Range("N64").Value = -10665

$ws.Range("H67").Value = 2831.5
$ws.Range("I67").Value = 1744.5
$ws.Range("J67").Value = 3375
$ws.Range("K67").Value = 5233.5
$ws.Range("L67").Value = 10125
$ws.Range("M67").Value = -4297.5
$ws.Range("N67").Value = -11997

$ws.Range("H92").Value = 350.75
$ws.Range("J92").Value = 367.66666
$ws.Range("L92").Value = 1102.99998
$ws.Range("N92").Value = -3598.99998

$ws.Range("H117").Value = 763.2
$ws.Range("I117").Value = 509
$ws.Range("J117").Value = 872.1429000000001
$ws.Range("K117").Value = 1527
$ws.Range("L117").Value = 2616.4287
$ws.Range("M117").Value = 1915
$ws.Range("N117").Value = -9500.4287

$ws.Range("H129").Value = 37697.3
$ws.Range("J129").Value = 41819.277
$ws.Range("L129").Value = 125457.831
$ws.Range("N129").Value = -135457.831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984

$ws.Range("H113").Value = 956.2308
$ws.Range("I113").Value = 684
$ws.Range("K113").Value = 684
$ws.Range("M113").Value = 1486

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1340
$ws.Range("I82").Value = 1475
$ws.Range("J82").Value = 800
$ws.Range("K82").Value = 1475
$ws.Range("L82").Value = 800
$ws.Range("M82").Value = -1114
$ws.Range("N82").Value = -1522

$ws.Range("H85").Value = 1340
$ws.Range("I85").Value = 1475
$ws.Range("J85").Value = 800
$ws.Range("K85").Value = 1475
$ws.Range("L85").Value = 800
$ws.Range("M85").Value = -227
$ws.Range("N85").Value = -3296

$ws.Range("H132").Value = 2996.12
$ws.Range("I132").Value = 2206.625
$ws.Range("K132").Value = 6619.875
$ws.Range("M132").Value = -4089.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 665.64703
$ws.Range("I107").Value = 531.1818
$ws.Range("J107").Value = 912.1667
$ws.Range("K107").Value = 1593.5454
$ws.Range("L107").Value = 2736.5001
$ws.Range("M107").Value = 326.4546
$ws.Range("N107").Value = -6576.5001

$ws.Range("H113").Value = 1149
$ws.Range("I113").Value = 698.6667
$ws.Range("K113").Value = 2096.0001
$ws.Range("M113").Value = 73.9998999999998

$ws.Range("H132").Value = 1448.9429
$ws.Range("I132").Value = 932.3077
$ws.Range("K132").Value = 2796.9231
$ws.Range("M132").Value = -266.9231

$ws.Range("H136").Value = 24157066
$ws.Range("I136").Value = 30866250
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 92598750
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -92596200
$ws.Range("N136").Value = -17100
